$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Semester 1 header row: extend the weekly date header into column F ---
$ws.Range("E27").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F27").Value = 43167

# --- Semester 1 task rows: add the new week's numbers ---
$ws.Range("F28").Value = 8
$ws.Range("F29").Value = 2

# --- "Add a semester" story block (rows 34-35) ---
$ws.Range("F34").Value = 15
$ws.Range("F35").Value = 5

# --- "Add classes to each semester" story block (rows 39-40) ---
$ws.Range("F39").Value = 15
$ws.Range("F40").Value = 5

# --- "Add course info / final grades" story block (rows 44-45) ---
$ws.Range("F44").Value = 15
$ws.Range("F45").Value = 5

# --- "Add grades / specify weight" story block (rows 49-50) ---
$ws.Range("F49").Value = 15
$ws.Range("F50").Value = 5

# --- Update the view: scroll down and select the "add a semester" block ---
$ws.Range("A34:A38").Select()
$excel.ActiveWindow.ScrollRow = 25
